# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G (header "K") values are recalculated from the new source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 3
    4  = 2
    5  = 2
    6  = 3
    7  = 3
    8  = 2
    9  = 9
    10 = 9
    11 = 7
    12 = 3
    13 = 8
    14 = 4
    15 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
